# Fill in the Timings schedule: replace template placeholder text with
# real meeting data (meeting IDs, course/venue names, dates, and the
# Tencent/Zoom "Way" column) for rows 2-10, and adjust formatting to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 -----------------------------------------------------------
# B2 used to hold placeholder text styled with the "微软雅黑" font (style 6);
# give it B3's plain style (style 2) before writing the real meeting id.
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B2").Value = 290526102
$ws.Range("C2").Value = ""
$ws.Range("F2").Value = "新时代"
$ws.Range("G2").Value = "Tencent"
# Row 2 had an oversized custom height (16.5); restore the default height.
$ws.Rows(2).AutoFit()

# --- Row 3 -------------------------------------------------------------
$ws.Range("B3").Value = 732938173
$ws.Range("F3").Value = "机电集成"
$ws.Range("G3").Value = "Tencent"

# --- Row 4 -------------------------------------------------------------
$ws.Range("B4").Value = 109147089
$ws.Range("F4").Value = "计算机通信"
$ws.Range("G4").Value = "Tencent"

# --- Row 5 -------------------------------------------------------------
$ws.Range("B5").Value = 298268214
$ws.Range("F5").Value = "机器视觉"
$ws.Range("G5").Value = "Tencent"

# --- Row 6 -------------------------------------------------------------
$ws.Range("A6").Value = 0.74722222222222223
$ws.Range("B6").Value = 403848798
$ws.Range("F6").Value = "先进制造技术与系统"
$ws.Range("G6").Value = "Tencent"

# --- Row 7 ---------------------------------------------------------
# Row 7 gains B/C cells matching the style used by rows 4-6 (style 5),
# loses its old D7 placeholder text, and gains a date in E7.
$ws.Range("B4:C4").Copy()
$ws.Range("B7:C7").PasteSpecial(-4122)
$ws.Range("A7").Value = 0.48680555555555555
$ws.Range("B7").Value = 124692830
$ws.Range("D7").ClearContents()
$ws.Range("E7").Value = 44679
$ws.Range("F7").Value = "院级高等讲堂"
$ws.Range("G7").Value = "Tencent"
$ws.Rows(7).RowHeight = 15

# --- Row 8 -------------------------------------------------------------
$ws.Range("A8").Value = 0.53472222222222221
$ws.Range("B8").Value = 84582727598
$ws.Range("E8").Value = 44684
$ws.Range("F8").Value = "院级高等讲堂"
$ws.Range("G8").Value = "Zoom"

# --- Row 9 -------------------------------------------------------------
$ws.Range("A9").Value = 0.79166666666666663
$ws.Range("B9").Value = 792644447
$ws.Range("E9").Value = 44689
$ws.Range("G9").Value = "Tencent"

# --- Row 10 --------------------------------------------------------
# Row 10 was essentially empty (only a styled E cell); build it out to
# match the same A/B/C/E/G pattern used by row 9, then fill in its data.
$ws.Range("A9:C9").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122)
$ws.Range("G9").Copy()
$ws.Range("G10").PasteSpecial(-4122)
$ws.Range("A10").Value = 0.75694444444444453
$ws.Range("B10").Value = 85621775445
$ws.Range("E10").Value = 44691
$ws.Range("G10").Value = "Zoom"
$ws.Rows(10).RowHeight = 15

# --- Selection -----------------------------------------------------
$null = $ws.Range("D15").Select()
